# Update "想去人数" (F column) counts on both the "展览" sheet and the
# "全部类型" sheet (which mirrors the same rows) to reflect refreshed
# scrape totals, per commit "Update gh-pages to output generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - source rows
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value  = 608
$wsExpo.Range("F5").Value  = 1145
$wsExpo.Range("F6").Value  = 14252
$wsExpo.Range("F7").Value  = 16178
$wsExpo.Range("F9").Value  = 81
$wsExpo.Range("F24").Value = 6486
$wsExpo.Range("F25").Value = 969
$wsExpo.Range("F26").Value = 10
$wsExpo.Range("F27").Value = 1111
$wsExpo.Range("F29").Value = 5680
$wsExpo.Range("F32").Value = 166
$wsExpo.Range("F33").Value = 4711

# Sheet "全部类型" (all types) - same events, shifted row numbers
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value  = 608
$wsAll.Range("F5").Value  = 1145
$wsAll.Range("F6").Value  = 14252
$wsAll.Range("F7").Value  = 16179
$wsAll.Range("F9").Value  = 81
$wsAll.Range("F25").Value = 6486
$wsAll.Range("F26").Value = 969
$wsAll.Range("F27").Value = 10
$wsAll.Range("F28").Value = 1111
$wsAll.Range("F31").Value = 5680
$wsAll.Range("F34").Value = 166
$wsAll.Range("F35").Value = 4711
